$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in July 2023 data on row 39
$ws.Range("C39").Value = 471
$ws.Range("D39").Value = 931
$ws.Range("E39").Value = 86.42
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Formula = "=(H39+I39)-(C39+D39+E39+F39+G39)"

# Update the active selection to J39
$ws.Range("J39").Select()
